$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 329; this shifts existing rows 329-336 down to 330-337.
$ws.Rows.Item(329).Insert()

# Populate the new row 329 with the weekly record.
$ws.Cells.Item(329, 1).Value = 6
$ws.Cells.Item(329, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(329, 3).Value = "Metropolitana"
$ws.Cells.Item(329, 4).Value = 44448
$ws.Cells.Item(329, 5).Value = 13
$ws.Cells.Item(329, 6).Value = 100112003
$ws.Cells.Item(329, 7).Value = "Ajo"
$ws.Cells.Item(329, 8).Value = "Chino"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 2550
$ws.Cells.Item(329, 11).Value = 14000
$ws.Cells.Item(329, 12).Value = 14500
$ws.Cells.Item(329, 13).Value = 14245
$ws.Cells.Item(329, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(329, 15).Value = "China"
$ws.Cells.Item(329, 16).Value = 1424
$ws.Cells.Item(329, 17).Value = 10
$ws.Cells.Item(329, 18).Value = "Hortaliza"
